# Updated symbol list on Fri Dec 23 22:58:37 UTC 2022 with GitHub Actions
# Refreshes the crypto "Price" (column D) and a couple of "Volume(1h)"
# (column E) values on Sheet1 to match the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    # Prefix with an apostrophe so a numeric-looking string (e.g. "246.14")
    # is stored as text, matching the sheet's existing inline-string cells
    # instead of being auto-converted into a number by Excel.
    $rng.Value = "'" + $NewValue
    # Restore the default "Normal" style so we don't leave a stray
    # quote-prefix / number-format flag on the cell.
    $rng.Style = "Normal"
}

# Column D ("Price") updates
Set-TextValue "D2"  "246.14"
Set-TextValue "D3"  "22.11"
Set-TextValue "D4"  "5.304"
Set-TextValue "D5"  "0.05873"
Set-TextValue "D6"  "3.382"
Set-TextValue "D7"  "6.382"
Set-TextValue "D8"  "0.8150"
Set-TextValue "D9"  "0.9570"
Set-TextValue "D11" "0.03603"
Set-TextValue "D12" "0.07321"
Set-TextValue "D13" "0.03051"
Set-TextValue "D14" "4.451"
Set-TextValue "D15" "0.09385"
Set-TextValue "D16" "0.001597"
Set-TextValue "D17" "0.04818"
Set-TextValue "D18" "0.0005902"
Set-TextValue "D19" "0.006206"
Set-TextValue "D20" "0.004082"
Set-TextValue "D21" "0.0009837"
Set-TextValue "D22" "0.00009703"
Set-TextValue "D23" "3.684"
Set-TextValue "D24" "2.181"
Set-TextValue "D25" "0.3267"
Set-TextValue "D27" "0.0002472"
Set-TextValue "D40" "0.03892"
Set-TextValue "D41" "0.006662"
Set-TextValue "D42" "0.1074"
Set-TextValue "D43" "0.003001"
Set-TextValue "D44" "0.005912"
Set-TextValue "D45" "0.00005662"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.7753"
Set-TextValue "D48" "0.08027"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.01010"

# Column E ("Volume(1h)") label updates
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E48").Value = "47BOLOBOLO"

Write-Host "Symbol list updated"
